# Updated cryptos list (price/volume refresh + a couple of row swaps
# for MXToken/VeChain and FraxShare/TheSandbox) as captured by the
# GitHub Actions commit on Sat May 20 02:03:47 UTC 2023.
#
# Column D ("Price") holds numeric-looking text (e.g. "308.92"), so each
# write is wrapped in a NumberFormat="@" / ClearFormats() dance to force
# Excel to keep it as a text cell instead of silently coercing it to a
# number - ClearFormats() afterwards restores the original (default)
# cell style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.906.71'
$c.ClearFormats()
$ws.Range("E2").Value = '  +0.27%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.813.31'
$c.ClearFormats()
$ws.Range("E3").Value = '  +0.66%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.07%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '308.92'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E6").Value = '  +0.04%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4654'
$c.ClearFormats()
$ws.Range("E7").Value = '  +0.85%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3666'
$c.ClearFormats()
$ws.Range("E8").Value = '  -1.27%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07352'
$c.ClearFormats()
$ws.Range("E9").Value = '  +0.19%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8672'
$c.ClearFormats()
$ws.Range("E10").Value = '  -0.05%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.28'
$c.ClearFormats()
$ws.Range("E11").Value = '  -0.43%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.815.68'
$c.ClearFormats()
$ws.Range("E12").Value = '  +2.76%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.373'
$c.ClearFormats()
$ws.Range("E13").Value = '  +0.47%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.07081'
$c.ClearFormats()
$ws.Range("E14").Value = '  +0.74%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.505'
$c.ClearFormats()
$ws.Range("E15").Value = '  -0.03%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '91.49'
$c.ClearFormats()
$ws.Range("E16").Value = '  +0.55%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.ClearFormats()
$ws.Range("E17").Value = '  +0.05%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008700'
$c.ClearFormats()
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("E19").Value = '  +0.05%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.66'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.00%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '26.930.81'
$c.ClearFormats()
$ws.Range("E21").Value = '  +0.32%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.301'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.29%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.61'
$c.ClearFormats()
$ws.Range("E23").Value = '  -0.45%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.043.57'
$c.ClearFormats()
$ws.Range("E24").Value = '  +2.01%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.892'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.50%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '150.16'
$c.ClearFormats()
$ws.Range("E26").Value = '  -0.77%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.165'
$c.ClearFormats()
$ws.Range("E27").Value = '  +1.89%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.28'
$c.ClearFormats()
$ws.Range("E28").Value = '  -0.28%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.265'
$c.ClearFormats()
$ws.Range("E29").Value = '  +0.61%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '115.40'
$c.ClearFormats()
$ws.Range("E30").Value = '  -0.05%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08938'
$c.ClearFormats()
$ws.Range("E31").Value = '  +0.66%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.7573'
$c.ClearFormats()
$ws.Range("E32").Value = '  +0.02%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.157'
$c.ClearFormats()
$ws.Range("E33").Value = '  +0.66%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.488'
$c.ClearFormats()
$ws.Range("E34").Value = '  +1.10%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.912'
$c.ClearFormats()
$ws.Range("E35").Value = '  -0.07%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  -1.83%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.05283'
$c.ClearFormats()
$ws.Range("E38").Value = '  +1.01%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.003'
$c.ClearFormats()
$ws.Range("E39").Value = '  +3.16%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.01952'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '7.193'
$c.ClearFormats()
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.5298'
$c.ClearFormats()
$ws.Range("E42").Value = '  +0.64%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.291'
$c.ClearFormats()
$ws.Range("E43").Value = '  -2.26%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.1652'
$c.ClearFormats()
$ws.Range("E44").Value = '  -0.06%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.387'
$c.ClearFormats()
$ws.Range("E45").Value = '  -0.99%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4861'
$c.ClearFormats()
$ws.Range("E46").Value = '  -2.86%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '10.41'
$c.ClearFormats()
$ws.Range("E47").Value = '  +1.16%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E48").Value = '  +0.07%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.660'
$c.ClearFormats()
$ws.Range("E49").Value = '  -0.02%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '102.77'
$c.ClearFormats()
$ws.Range("E50").Value = '  -0.78%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.06286'
$c.ClearFormats()
$ws.Range("E51").Value = '  -0.06%  '
